# DEG Number of Identified Genes.xlsx - update
#
# 1. Rename the (only) worksheet from
#      "Num Stat Sig Genes Identified"
#    to
#      "M_MUT_and_WT_M_P30_CORT"
# 2. Move the active selection on that sheet from A11:K12 to the single
#    cell C24.
# 3. Reposition the workbook window on screen (xWindow goes from -108 to
#    -23148; yWindow stays -108) to match the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet.
$ws.Name = "M_MUT_and_WT_M_P30_CORT"

# 2) Update the saved selection/active cell.
$ws.Range("C24").Select()

# 3) Update the saved window position.
$win = $excel.ActiveWindow
$win.Left = -23148
$win.Top = -108
